$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaLabel = "Meta description"
$metaRest  = ": Discover the gameplay, bonuses, graphics and winning potential of Attila slot. Play it for free and win big!"
$metaPara.Range.Text = $metaLabel + $metaRest

$metaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $metaLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the duplicated "Play Attila Slot for Free - Review" (bold)
#    paragraph that sits near the end of the document, right before the
#    italic meta-description-style paragraph.
#    (Paragraph.Range.Text carries a trailing paragraph-mark character, so
#    trim it off before comparing.)
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Play Attila Slot for Free - Review" -and $i -gt 1) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph (the old meta
#    description) with the new image-generation prompt, keeping its
#    italic run formatting intact.
# ---------------------------------------------------------------------------
$oldText = "Discover the gameplay, bonuses, graphics and winning potential of Attila slot. Play it for free and win big!"
$newText = 'Please create a cartoon-style feature image for the game "Attila". The image should feature a happy Maya warrior wearing glasses. The image should be colorful and eye-catching, and should include elements from the game such as arrows, axes, and the Hunnic army. The image should be dynamic and convey the excitement and adventure of the game. Please ensure that the image is high-quality and visually appealing, and that it accurately represents the key features and functionalities of the game.'

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq $oldText) {
        $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
        $rng.Text = $newText
        break
    }
}
